$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.34"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.00"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.338"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05959"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8104"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9641"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1426"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07391"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03400"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03053"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09403"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.995"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001591"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04806"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005913"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006238"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005137"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009830"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00009706"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.744"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.186"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03914"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1074"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002712"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.006534"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005820"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005317"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8505"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03526"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01011"
